# 70827 Ultrakatty with Warrior Lucy - 300% BOM update
# "Added tiles and plates, updated 70827 BOM"

$wb = $excel.ActiveWorkbook

$fmtCurrency = "_-* #,##0.00\ ""kr""_-;\-* #,##0.00\ ""kr""_-;_-* ""-""??\ ""kr""_-;_-@_-"
$fmtTimeAmPm = "[$-F400]h:mm:ss\ AM/PM"
$fmtTime = "h:mm:ss"

# ---------------------------------------------------------------
# Sheet "Black" (1) - add Act time (H2) for the first part
# ---------------------------------------------------------------
$wsBlack = $wb.Worksheets.Item(1)

$wsBlack.Range("H2").Value = 0.027083333333333334
$wsBlack.Range("H2").NumberFormat = $fmtTimeAmPm

$wsBlack.Range("A1:XFD1").Select()

# ---------------------------------------------------------------
# Sheet "Red" (3) - fill in filament / cost / time columns
# ---------------------------------------------------------------
$wsRed = $wb.Worksheets.Item(3)

$wsRed.Range("D5").Value = 21.43
$wsRed.Range("E5").Value = 64.42
$wsRed.Range("F5").Value = 19.260000000000002
$wsRed.Range("G5").Value = 0.23608796296296297
$wsRed.Range("G5").NumberFormat = $fmtTime

$wsRed.Range("D7").Value = 4.93
$wsRed.Range("E7").Value = 14.81
$wsRed.Range("F7").Value = 4.43
$wsRed.Range("G7").Value = 0.05950231481481482
$wsRed.Range("G7").NumberFormat = $fmtTime

$wsRed.Range("D8").Value = 35.090000000000003
$wsRed.Range("E8").Value = 105.49
$wsRed.Range("F8").Value = 31.54
$wsRed.Range("G8").Value = 0.38418981481481485
$wsRed.Range("G8").NumberFormat = $fmtTime

$wsRed.Range("D9").Value = 30.53
$wsRed.Range("E9").Value = 91.79
$wsRed.Range("F9").Value = 27.45
$wsRed.Range("G9").Value = 0.28556712962962966
$wsRed.Range("G9").NumberFormat = $fmtTime

$wsRed.Range("D10").Value = 5.8
$wsRed.Range("E10").Value = 17.43
$wsRed.Range("F10").Value = 5.21
$wsRed.Range("G10").Value = 0.05451388888888889
$wsRed.Range("G10").NumberFormat = $fmtTime

$wsRed.Range("D11").Value = 13.22
$wsRed.Range("E11").Value = 39.74
$wsRed.Range("F11").Value = 11.88
$wsRed.Range("G11").Value = 0.13328703703703704
$wsRed.Range("G11").NumberFormat = $fmtTime

$wsRed.Range("D12").Value = 18.47
$wsRed.Range("E12").Value = 55.52
$wsRed.Range("F12").Value = 16.600000000000001
$wsRed.Range("G12").Value = 0.19582175925925926
$wsRed.Range("G12").NumberFormat = $fmtTime

$wsRed.Range("D14").Value = 25.88
$wsRed.Range("E14").Value = 77.81
$wsRed.Range("F14").Value = 23.26
$wsRed.Range("G14").Value = 0.25655092592592593
$wsRed.Range("G14").NumberFormat = $fmtTime

$wsRed.Range("E6").Select()

# ---------------------------------------------------------------
# Sheet "White" (6) - add header row + filament/cost/time columns
# ---------------------------------------------------------------
$wsWhite = $wb.Worksheets.Item(6)

$wsWhite.Rows.Item(1).Insert()

$wsWhite.Range("A1").Value = "Item No"
$wsWhite.Range("B1").Value = "Item"
$wsWhite.Range("C1").Value = "Copies"
$wsWhite.Range("D1").Value = "Filament (m)"
$wsWhite.Range("E1").Value = "Filament (g)"
$wsWhite.Range("F1").Value = "Cost (kr)"
$wsWhite.Range("F1").NumberFormat = $fmtCurrency
$wsWhite.Range("G1").Value = "Est time"
$wsWhite.Range("G1").NumberFormat = $fmtTimeAmPm
$wsWhite.Range("H1").Value = "Act time"
$wsWhite.Range("H1").NumberFormat = $fmtTimeAmPm

$wsWhite.Range("D2").Value = 1.96
$wsWhite.Range("E2").Value = 5.89
$wsWhite.Range("F2").Value = 1.76
$wsWhite.Range("F2").NumberFormat = $fmtCurrency
$wsWhite.Range("G2").Value = 0.020879629629629626
$wsWhite.Range("G2").NumberFormat = $fmtTime
$wsWhite.Range("H2").Value = 0.022222222222222223
$wsWhite.Range("H2").NumberFormat = $fmtTime

$wsWhite.Range("H3").Select()

# ---------------------------------------------------------------
# Sheet "Blue" (2) - add header row + filament/cost/time/notes
# (this is the last sheet touched, so it stays the active tab)
# ---------------------------------------------------------------
$wsBlue = $wb.Worksheets.Item(2)

$wsBlue.Rows.Item(1).Insert()

$wsBlue.Range("A1").Value = "Item No"
$wsBlue.Range("B1").Value = "Item"
$wsBlue.Range("C1").Value = "Copies"
$wsBlue.Range("D1").Value = "Filament (m)"
$wsBlue.Range("E1").Value = "Filament (g)"
$wsBlue.Range("F1").Value = "Cost (kr)"
$wsBlue.Range("F1").NumberFormat = $fmtCurrency
$wsBlue.Range("G1").Value = "Est time"
$wsBlue.Range("G1").NumberFormat = $fmtTimeAmPm
$wsBlue.Range("H1").Value = "Act time"
$wsBlue.Range("H1").NumberFormat = $fmtTimeAmPm
$wsBlue.Range("I1").Value = "Notes"

$wsBlue.Range("D3").Value = 3.92
$wsBlue.Range("E3").Value = 11.7
$wsBlue.Range("F3").Value = 3.5
$wsBlue.Range("F3").NumberFormat = $fmtCurrency
$wsBlue.Range("G3").Value = 0.038495370370370367
$wsBlue.Range("G3").NumberFormat = $fmtTime

$wsBlue.Range("D4").Value = 28.92
$wsBlue.Range("E4").Value = 86.24
$wsBlue.Range("F4").Value = 25.79
$wsBlue.Range("F4").NumberFormat = $fmtCurrency
$wsBlue.Range("G4").Value = 0.31452546296296297
$wsBlue.Range("G4").NumberFormat = $fmtTime

$wsBlue.Range("D5").Value = 7.65
$wsBlue.Range("E5").Value = 22.81
$wsBlue.Range("F5").Value = 6.82
$wsBlue.Range("F5").NumberFormat = $fmtCurrency
$wsBlue.Range("G5").Value = 0.070717592592592596
$wsBlue.Range("G5").NumberFormat = $fmtTime
$wsBlue.Range("H5").Value = 0.071527777777777787
$wsBlue.Range("H5").NumberFormat = $fmtTime

$wsBlue.Range("D6").Value = 14.84
$wsBlue.Range("E6").Value = 44.27
$wsBlue.Range("F6").Value = 13.24
$wsBlue.Range("F6").NumberFormat = $fmtCurrency
$wsBlue.Range("G6").Value = 0.13685185185185186
$wsBlue.Range("G6").NumberFormat = $fmtTime
$wsBlue.Range("H6").Value = 0.1388888888888889
$wsBlue.Range("H6").NumberFormat = $fmtTime

$wsBlue.Range("D7").Value = 1.08
$wsBlue.Range("E7").Value = 3.21
$wsBlue.Range("F7").Value = 0.96
$wsBlue.Range("F7").NumberFormat = $fmtCurrency
$wsBlue.Range("G7").Value = 0.014606481481481482
$wsBlue.Range("G7").NumberFormat = $fmtTime
$wsBlue.Range("I7").Value = "Cut @ 7,8"

$wsBlue.Range("D8").Value = 3.98
$wsBlue.Range("E8").Value = 11.6
$wsBlue.Range("F8").Value = 3.47
$wsBlue.Range("F8").NumberFormat = $fmtCurrency
$wsBlue.Range("G8").Value = 0.038657407407407404
$wsBlue.Range("G8").NumberFormat = $fmtTime
$wsBlue.Range("I8").Value = "45 deg support for groove"

$wsBlue.Activate()
$wsBlue.Range("I9").Select()

Write-Host "Edit complete"
